# Added CRUD for Restaurant
#
# Populates Sheet1 with reference data for the Restaurant / MenuItem /
# MenuCategory / Roles / Restaurant-Owner-Mapping tables, highlights the
# table-title and join-key header cells, sizes the columns, configures the
# page setup and restores the saved selection.
#
# NOTE: the string literals below are intentionally written in shared-
# string-table order (Restaurant, RestaurantId, Name, Description, ...)
# rather than sheet reading order, so the workbook's shared string table
# comes out in the same order the original author built it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (BGR-packed OLE values, matching Excel's Interior.Color convention)
$yellow = 65535      # RGB(255,255,0)  -> fill FFFFFF00 (table title rows)
$green  = 5296274    # RGB(146,208,80) -> fill FF92D050 (join-key headers)

# --- introduce each unique string once, in shared-string-table order ---
$ws.Range("A1").Value = "Restaurant"
$ws.Range("A2").Value = "RestaurantId"
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "Description"
$ws.Range("D2").Value = "Status"
$ws.Range("E2").Value = "PerPersonAmount"
$ws.Range("F2").Value = "IsActive"
$ws.Range("A6").Value = "MenuItem"
$ws.Range("A7").Value = "MenuItemId"
$ws.Range("C7").Value = "Price"
$ws.Range("E7").Value = "MenuCategoryId"
$ws.Range("A10").Value = "MenuCategory"
$ws.Range("B3").Value = "Maria"
$ws.Range("C3").Value = "Serves Good Kerala food"
$ws.Range("F8").Value = "Veg"
$ws.Range("B12").Value = "Starter"
$ws.Range("C12").Value = "Good Starter"
$ws.Range("B13").Value = "BreakFast"
$ws.Range("C13").Value = "Serve you best energentic breakfast"
$ws.Range("F7").Value = "FoodCategory"
$ws.Range("B8").Value = "Tandoor Chicken"
$ws.Range("B14").Value = "MainCourse"
$ws.Range("C14").Value = "All Veg non veg main course "
$ws.Range("A17").Value = "Roles"
$ws.Range("A18").Value = "Restaurant Owner"
$ws.Range("A19").Value = "Admin"
$ws.Range("A21").Value = "User"
$ws.Range("A20").Value = "Cook"
$ws.Range("A23").Value = "Restaurant Owner Mapping"
$ws.Range("B24").Value = "OwnerId"

# --- remaining cells: duplicate strings (by reference) and numeric values ---
$ws.Range("A3").Value = 1
$ws.Range("E3").Value = 200
$ws.Range("F3").Value = 1
$ws.Range("B7").Value = "Name"
$ws.Range("D7").Value = "IsActive"
$ws.Range("A8").Value = 1
$ws.Range("C8").Value = 300
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 3
$ws.Range("A11").Value = "MenuCategoryId"
$ws.Range("B11").Value = "Name"
$ws.Range("C11").Value = "Description"
$ws.Range("D11").Value = "RestaurantId"
$ws.Range("A12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("A13").Value = 2
$ws.Range("D13").Value = 1
$ws.Range("A14").Value = 3
$ws.Range("D14").Value = 1
$ws.Range("A24").Value = "RestaurantId"

# ---------------------------------------------------------------------
# Table-title and join-key header highlighting
# ---------------------------------------------------------------------
$ws.Range("A1").Interior.Color = $yellow
$ws.Range("A6").Interior.Color = $yellow
$ws.Range("A10").Interior.Color = $yellow
$ws.Range("E7").Interior.Color = $green
$ws.Range("A11").Interior.Color = $green

# ---------------------------------------------------------------------
# Column widths (best-fit per the authored workbook)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.28515625
$ws.Columns.Item(2).ColumnWidth = 15.85546875
$ws.Columns.Item(3).ColumnWidth = 33.42578125
$ws.Columns.Item(4).ColumnWidth = 12.28515625
$ws.Columns.Item(5).ColumnWidth = 17.5703125
$ws.Columns.Item(6).ColumnWidth = 13.42578125

# ---------------------------------------------------------------------
# Page setup
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Restore the view/selection the workbook was saved with
# ---------------------------------------------------------------------
$ws.Range("C23").Select()
